$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header row (row 1) ---
$ws.Range("A1").Value = "论文标题"
$ws.Range("B1").Value = "年份"
$ws.Range("C1").Value = "类型"
$ws.Range("D1").Value = "解决问题的方法"
$ws.Range("E1").Value = "实验语料"
$ws.Range("F1").Value = "缺点"
$ws.Range("G1").Value = "可以抄的部分"

# --- Data row (row 2) ---
# Values are assigned in the same order the original author typed them
# (this also drives shared-string insertion order on save).
$ws.Range("A2").Value = "A Chinese Dataset with Negative Full Forms for General Abbreviation`nPrediction`n"
$ws.Range("C2").Value = "期刊"
$ws.Range("G2").Value = "Abbreviation processing mainly consists of three tasks, that`nis, abbreviation expansion, abbreviation recognition, and`nabbreviation prediction. Expanding the short form of an`nexpression to its full form is called abbreviation expansion.`nExtracting the short form and full form pairs from the context is called abbreviation recognition. Abbreviation prediction refers to predicting the short form of an expression according to its full form.`nSun et al. (2009) shows that better abbreviation prediction will help improve the performance of abbreviation recognition.`nRobust approach to abbreviating terms: A discriminative latent`nvariable model with global information."
$ws.Range("E2").Value = "人民日报语料库"
$ws.Range("D2").Value = "用不同的算法(CRF、LSTM）生成了`n包含NFF(没有简称的全称）的数据集。其实没什么贡献啊，这也能发论文啊，还是北大的"
$ws.Range("B2").Value = 2017

# --- Wrap text for the cells that hold long text ---
$ws.Range("A2").WrapText = $true
$ws.Range("D2").WrapText = $true
$ws.Range("G2").WrapText = $true

# --- Column widths (closest achievable values to the target stored widths) ---
$ws.Range("A1").EntireColumn.ColumnWidth = 24.3
$ws.Range("D1").EntireColumn.ColumnWidth = 29.5
$ws.Range("E1").EntireColumn.ColumnWidth = 24.15
$ws.Range("G1").EntireColumn.ColumnWidth = 57

# --- Row height for data row ---
$ws.Range("A2").EntireRow.RowHeight = 101.25

# --- Default row height for the sheet ---
$ws.StandardHeight = 14.25

# --- Selection state ---
[void]$ws.Range("G17").Select()
